$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the Price (column D) and Volume(1h) (column E) figures for the
# coin rows that changed in the latest symbol-list scrape.
#
# The source cells are plain text (inline strings), e.g. "324.79" or "2.42%,"
# not numbers. Writing the new value with a leading apostrophe forces Excel to
# keep it as literal text instead of coercing it into a number/percentage; the
# ClearFormats() call right after removes the "quote prefix" marker style that
# the apostrophe otherwise leaves behind, so the cell ends up with the same
# (unstyled) look as before the edit.

$updates = @{
    "D2" = "324.79"
    "E2" = "2.42%"
    "D3" = "39.96"
    "E3" = "5.64%"
    "D4" = "5.870"
    "E4" = "13.37%"
    "D5" = "0.07993"
    "E5" = "-0.50%"
    "D6" = "4.583"
    "E6" = "2.05%"
    "D7" = "8.707"
    "E7" = "2.31%"
    "D8" = "1.912"
    "E8" = "-0.84%"
    "D9" = "2.942"
    "E9" = "-1.94%"
    "D10" = "0.9392"
    "E10" = "0.08%"
    "D11" = "0.1258"
    "E11" = "-2.53%"
    "D12" = "0.1961"
    "E12" = "1.14%"
    "D13" = "8.821"
    "E13" = "33.65%"
    "D14" = "0.09153"
    "E14" = "1.19%"
    "D15" = "0.03571"
    "E15" = "4.89%"
    "D16" = "0.09615"
    "E16" = "0.85%"
    "D17" = "0.001302"
    "E17" = "-6.81%"
    "D18" = "0.006169"
    "E18" = "-1.01%"
    "D19" = "3.347"
    "E19" = "-0.92%"
    "E20" = "0.12%"
    "D21" = "0.1434"
    "E21" = "8.97%"
    "E22" = "-0.39%"
    "D23" = "0.04438"
    "E23" = "0.95%"
    "D24" = "0.001262"
    "E24" = "2.58%"
    "D25" = "0.004347"
    "E25" = "1.78%"
    "E26" = "-13.89%"
    "E27" = "0.13%"
    "D39" = "0.02421"
    "E39" = "2.99%"
    "D40" = "0.05254"
    "E40" = "1.99%"
    "D41" = "0.007480"
    "E41" = "-1.89%"
    "D42" = "0.1411"
    "E42" = "0.67%"
    "D43" = "0.008697"
    "E43" = "0.13%"
    "D44" = "0.002105"
    "E44" = "-0.17%"
    "D45" = "0.01055"
    "E45" = "19.32%"
    "D46" = "0.00006871"
    "E46" = "6.51%"
    "D47" = "0.00000000753"
    "E47" = "0.46%"
    "D48" = "0.002885"
    "D49" = "0.001426"
    "E49" = "-15.56%"
    "D50" = "0.00002109"
    "E50" = "0.46%"
    "D51" = "0.0002008"
    "E51" = "0.46%"
}

foreach ($addr in $updates.Keys) {
    $range = $ws.Range($addr)
    $range.Value = "'" + $updates[$addr]
    $range.ClearFormats()
}
